# Applies the v2.5 interface update: clears the stray empty F47/G47 cells
# and appends the new incident rows (48-64) recorded on 2024-05-31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F47 and G47 were placeholder empty inline-string cells with no data;
# remove them entirely so the row ends at column E/H like the rest.
$ws.Range("F47:G47").ClearContents()

# Row 48
$ws.Range("A48").Value = "WC48 P5F"
$ws.Range("B48").Value = "Cámara no detecta skeleton"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "2024-05-31"
$ws.Range("D48").Value = "11:27:23"
$ws.Range("E48").Value = "Mañana"
$ws.Range("F48").Value = "11:27:28"
$ws.Range("G48").Value = "0:00:05"
$ws.Range("H48").Value = "N/A"

# Row 49
$ws.Range("A49").Value = "WC48 P5F"
$ws.Range("B49").Value = "Cámara no detecta foams"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "2024-05-31"
$ws.Range("D49").Value = "11:27:39"
$ws.Range("E49").Value = "Mañana"
$ws.Range("F49").Value = "11:27:41"
$ws.Range("G49").Value = "0:00:02"
$ws.Range("H49").Value = "0.12 minutos"

# Row 50
$ws.Range("A50").Value = "WC48 P5F"
$ws.Range("B50").Value = "Etiquetadora"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "2024-05-31"
$ws.Range("D50").Value = "11:31:03"
$ws.Range("E50").Value = "Mañana"
$ws.Range("F50").Value = "11:31:06"
$ws.Range("G50").Value = "0:00:03"
$ws.Range("H50").Value = "0.19 minutos"

# Row 51
$ws.Range("A51").Value = "WC48 P5F"
$ws.Range("B51").Value = "Etiquetadora"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "2024-05-31"
$ws.Range("D51").Value = "11:31:05"
$ws.Range("E51").Value = "Mañana"
$ws.Range("F51").Value = "11:31:07"
$ws.Range("G51").Value = "0:00:02"
$ws.Range("H51").Value = "1.26 minutos"

# Row 52
$ws.Range("A52").Value = "WC48 P5F"
$ws.Range("B52").Value = "Cámara no detecta Power CP"
$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "2024-05-31"
$ws.Range("D52").Value = "11:31:11"
$ws.Range("E52").Value = "Mañana"
$ws.Range("F52").Value = "11:31:14"
$ws.Range("G52").Value = "0:00:03"
$ws.Range("H52").Value = "0.95 minutos"

# Row 53
$ws.Range("A53").Value = "WC48 P5F"
$ws.Range("B53").Value = "Cámara no detecta Power CP"
$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "2024-05-31"
$ws.Range("D53").Value = "11:31:12"
$ws.Range("E53").Value = "Mañana"
$ws.Range("F53").Value = "11:31:14"
$ws.Range("G53").Value = "0:00:02"
$ws.Range("H53").Value = "0.78 minutos"

# Row 54
$ws.Range("A54").Value = "WC48 P5F"
$ws.Range("B54").Value = "Cámara no detecta busbar"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "2024-05-31"
$ws.Range("D54").Value = "12:09:27"
$ws.Range("E54").Value = "Mañana"
$ws.Range("F54").Value = "12:09:33"
$ws.Range("G54").Value = "0:00:06"
$ws.Range("H54").Value = "N/A"

# Row 55
$ws.Range("A55").Value = "WC48 P5F"
$ws.Range("B55").Value = "AOI (malla)"
$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value = "2024-05-31"
$ws.Range("D55").Value = "12:09:32"
$ws.Range("E55").Value = "Mañana"
$ws.Range("F55").Value = "12:09:34"
$ws.Range("G55").Value = "0:00:02"
$ws.Range("H55").Value = "N/A"

# Row 56
$ws.Range("A56").Value = "WC48 P5F"
$ws.Range("B56").Value = "AOI (malla)"
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "2024-05-31"
$ws.Range("D56").Value = "12:09:37"
$ws.Range("E56").Value = "Mañana"
$ws.Range("F56").Value = "12:09:38"
$ws.Range("G56").Value = "0:00:01"
$ws.Range("H56").Value = "0.09 minutos"

# Row 57
$ws.Range("A57").Value = "WC48 P5F"
$ws.Range("B57").Value = "Cámara no detecta skeleton"
$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "2024-05-31"
$ws.Range("D57").Value = "12:14:01"
$ws.Range("E57").Value = "Mañana"
$ws.Range("F57").Value = "12:14:03"
$ws.Range("G57").Value = "0:00:02"
$ws.Range("H57").Value = "0.09 minutos"

# Row 58
$ws.Range("A58").Value = "WC48 P5F"
$ws.Range("B58").Value = "AOI (malla)"
$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "2024-05-31"
$ws.Range("D58").Value = "12:14:14"
$ws.Range("E58").Value = "Mañana"
$ws.Range("F58").Value = "12:14:19"
$ws.Range("G58").Value = "0:00:05"
$ws.Range("H58").Value = "1.52 minutos"

# Row 59
$ws.Range("A59").Value = "WV50 FILTER"
$ws.Range("B59").Value = "No coloca bien el core"
$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "2024-05-31"
$ws.Range("D59").Value = "12:23:32"
$ws.Range("E59").Value = "Mañana"
$ws.Range("F59").Value = "12:23:38"
$ws.Range("G59").Value = "0:00:06"
$ws.Range("H59").Value = "N/A"

# Row 60
$ws.Range("A60").Value = "WV50 FILTER"
$ws.Range("B60").Value = "Fallo visión core"
$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "2024-05-31"
$ws.Range("D60").Value = "12:23:37"
$ws.Range("E60").Value = "Mañana"
$ws.Range("F60").Value = "12:23:39"
$ws.Range("G60").Value = "0:00:02"
$ws.Range("H60").Value = "N/A"

# Row 61
$ws.Range("A61").Value = "SPL"
$ws.Range("B61").Value = "Fallo dispensación glue"
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "2024-05-31"
$ws.Range("D61").Value = "12:30:52"
$ws.Range("E61").Value = "Mañana"
$ws.Range("H61").Value = "N/A"

# Row 62
$ws.Range("A62").Value = "SPL"
$ws.Range("B62").Value = "Soldadura defectuosa"
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").Value = "2024-05-31"
$ws.Range("D62").Value = "12:30:54"
$ws.Range("E62").Value = "Mañana"
$ws.Range("H62").Value = "N/A"

# Row 63
$ws.Range("A63").Value = "SPL"
$ws.Range("B63").Value = "Soldadura defectuosa"
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "2024-05-31"
$ws.Range("D63").Value = "12:30:55"
$ws.Range("E63").Value = "Mañana"
$ws.Range("H63").Value = "0.03 minutos"

# Row 64
$ws.Range("A64").Value = "WV50 FILTER"
$ws.Range("B64").Value = "Traza"
$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "2024-05-31"
$ws.Range("D64").Value = "12:48:55"
$ws.Range("E64").Value = "Mañana"
$ws.Range("F64").Value = "12:48:57"
$ws.Range("G64").Value = "0:00:02"
$ws.Range("H64").Value = "N/A"

